# Auto-generated edit script: updates cryptos list values
# per commit 'Updated cryptos list on Sat Aug 26 06:53:45 UTC 2023 with GitHub Actions'
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.098.92"
$ws.Range("E2").Value = "  -0.06%  "

# Row 3
$ws.Range("D3").Value = "1.654.81"
$ws.Range("E3").Value = "  -0.11%  "

# Row 4
$ws.Range("E4").Value = "  -0.27%  "

# Row 5
$ws.Range("D5").Value = "'217.62"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.66%  "

# Row 6
$ws.Range("D6").Value = "'0.5256"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.64%  "

# Row 7
$ws.Range("E7").Value = "  -0.22%  "

# Row 8
$ws.Range("E8").Value = "  -0.76%  "

# Row 9
$ws.Range("D9").Value = "'0.06355"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.51%  "

# Row 10
$ws.Range("D10").Value = "'20.39"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.66%  "

# Row 11
$ws.Range("D11").Value = "'0.07805"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.11%  "

# Row 12
$ws.Range("D12").Value = "'4.511"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.04%  "

# Row 13
$ws.Range("D13").Value = "1.663.31"
$ws.Range("E13").Value = "  +0.23%  "

# Row 14
$ws.Range("D14").Value = "'0.5479"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.27%  "

# Row 15
$ws.Range("D15").Value = "0.0₅8212"
$ws.Range("E15").Value = "  +1.28%  "

# Row 16
$ws.Range("D16").Value = "'65.41"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.13%  "

# Row 17
$ws.Range("D17").Value = "26.128.06"
$ws.Range("E17").Value = "  -0.05%  "

# Row 18
$ws.Range("E18").Value = "  -0.29%  "

# Row 19
$ws.Range("D19").Value = "'4.590"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.46%  "

# Row 20
$ws.Range("D20").Value = "'191.51"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.20%  "

# Row 21
$ws.Range("E21").Value = "  -0.01%  "

# Row 22
$ws.Range("D22").Value = "'6.021"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.02%  "

# Row 23
$ws.Range("E23").Value = "  -0.29%  "

# Row 24
$ws.Range("D24").Value = "'142.08"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.59%  "

# Row 25
$ws.Range("E25").Value = "  +1.31%  "

# Row 26
$ws.Range("D26").Value = "'7.254"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.26%  "

# Row 27
$ws.Range("D27").Value = "'16.10"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.23%  "

# Row 28
$ws.Range("D28").Value = "'1.429"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.97%  "

# Row 29
$ws.Range("D29").Value = "'0.05894"
$ws.Range("D29").Style = "Normal"

# Row 30
$ws.Range("D30").Value = "'1.274"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.55%  "

# Row 31
$ws.Range("D31").Value = "'3.515"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.67%  "

# Row 32
$ws.Range("D32").Value = "'3.257"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.12%  "

# Row 33
$ws.Range("D33").Value = "'1.594"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.46%  "

# Row 34
$ws.Range("D34").Value = "'0.9509"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.40%  "

# Row 35
$ws.Range("D35").Value = "'2.785"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.40%  "

# Row 36
$ws.Range("D36").Value = "'2.412"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.61%  "

# Row 37
$ws.Range("D37").Value = "'0.5698"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.51%  "

# Row 38
$ws.Range("D38").Value = "'0.01619"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.76%  "

# Row 39
$ws.Range("D39").Value = "'5.814"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.46%  "

# Row 40
$ws.Range("D40").Value = "'0.8479"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.80%  "

# Row 41
$ws.Range("E41").Value = "  -0.16%  "

# Row 42
$ws.Range("D42").Value = "1.030.09"
$ws.Range("E42").Value = "  +2.09%  "

# Row 43
$ws.Range("D43").Value = "'102.84"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.44%  "

# Row 44
$ws.Range("D44").Value = "1.798.75"
$ws.Range("E44").Value = "  +0.06%  "

# Row 45
$ws.Range("D45").Value = "'57.13"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.99%  "

# Row 46
$ws.Range("E46").Value = "  -0.32%  "

# Row 47
$ws.Range("D47").Value = "'0.4304"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.56%  "

# Row 48
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "'1.476"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.85%  "

# Row 49
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "'0.05154"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.28%  "

# Row 50
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'7.840"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.12%  "

# Row 51
$ws.Range("D51").Value = "'0.09694"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.12%  "

